$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the header row text (TAC-3844: fix edit texts in import trips
#    excel files). Each header cell gets an English line + an Arabic line
#    separated by a line break.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Trip Reference No *  `nالرقم المرجعي للرحلة"
$ws.Range("B1").Value = "Pick-up date start  (yyyy/MM/dd)`nتاريخ بداية التحميل "
$ws.Range("C1").Value = "Pick-up date end  (yyyy/MM/dd)`nتاريخ نهاية التحميل "
$ws.Range("D1").Value = "Approximate total value of goods`nالقيمة التقريبة للبضاعة "
$ws.Range("E1").Value = "Shipper notes`nملاحظات الشاحن"
$ws.Range("F1").Value = "Needs delivery note?`nبحاجة إلى مذكرة تسليم؟"
$ws.Range("G1").Value = "Has attachment?`nهل يوجد مرفقات ؟"
$ws.Range("H1").Value = "Pick-up facility*`nمنشأة التحميل"
$ws.Range("I1").Value = "Drop-off facility*`nمنشأة التنزيل"
$ws.Range("J1").Value = "Sender name*`nاسم المرسل "
$ws.Range("K1").Value = "Reciver name*`nاسم المستلم "

# ---------------------------------------------------------------------------
# 2. Re-style the template: thin black box borders around the header and the
#    input columns (B, C), smaller 8pt fonts, centred/wrapped header, and a
#    text-formatted, bordered entry area for the first 29 data rows.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:K1")
$headerRange.Font.Name = "Open Sans"
$headerRange.Font.Size = 8
$headerRange.Font.Color = 0
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108
$headerRange.WrapText = $true
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Color = 0
$ws.Rows(1).RowHeight = 25.5

$bodyRange = $ws.Range("B2:C30")
$bodyRange.Font.Name = "Calibri"
$bodyRange.Font.Size = 8
$bodyRange.NumberFormat = "@"
$bodyRange.Borders.LineStyle = 1
$bodyRange.Borders.Color = 0

# First entry cell keeps an explicit left alignment (matches the rest of the
# text cells visually, as produced by the source template).
$ws.Range("B2").HorizontalAlignment = -4131

# ---------------------------------------------------------------------------
# 3. Data validation: the Yes/No list now only spans columns F:G (previously
#    "G1:G1048576 F1:F1048576", now the single contiguous "F1:G1048576").
# ---------------------------------------------------------------------------
$ws.Range("F1:G1048576").Validation.Delete()
$ws.Range("F1:G1048576").Validation.Add(3, 1, 1, "yes,no")
$ws.Range("F1:G1048576").Validation.IgnoreBlank = $true
$ws.Range("F1:G1048576").Validation.InCellDropdown = $true
$ws.Range("F1:G1048576").Validation.ShowInput = $true
$ws.Range("F1:G1048576").Validation.ShowError = $true

# ---------------------------------------------------------------------------
# 4. Sheet view: scroll so row 20 is at the top and select A1:XFD30.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 20
$ws.Range("A1:XFD30").Select()
